$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column at C ("Request Status" moves from the end to position C) ---
$ws.Columns("C").Insert()

# Give the new column C the same width as column B (raw OOXML width 15).
$ws.Columns("C").ColumnWidth = 14.17

# Fill header of the new column with the "Request Status" label.
# (After the insert, the old trailing "Request Status" header has shifted from
#  AG1 to AH1, so AH1 already holds that text - overwrite C1 with a fresh copy
#  of the same label rather than deleting/re-adding a column.)
$ws.Range("C1").Value = "Request Status"

# --- Turn the old trailing header (now at AH1, still "Request Status") into the
#     brand-new "Asset Status" column, keeping its existing header styling ---
$ws.Range("AH1").Value = "Asset Status"

# --- Refresh the AutoFilter so it spans the new range A1:AH1 ---
$ws.AutoFilterMode = $false
$ws.Range("A1:AH1").AutoFilter()

# --- Keep the _FilterDatabase defined name in sync with the new range ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Data!_FilterDatabase") {
        $n.RefersTo = "=Data!`$A`$1:`$AH`$1"
    }
}

# --- Update the active selection shown in the sheet view ---
$ws.Range("C2").Select()
